$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 <- old row 4
$ws.Range("A3").Value = 87535081
$ws.Range("B3").Value = 103178
$ws.Range("E3").Value = 221141
$ws.Range("Q3").Value = 693869.236562217
$ws.Range("R3").Value = 6685521.203787691
$ws.Range("F3").Value = "Gullviva"
$ws.Range("G3").Value = "Primula veris"
$ws.Range("J3").Value = "plantor/tuvor"
$ws.Range("AC3").Value = "Förekomst längs väg nr 1120 S om vägbanan i innerslänt."
$ws.Range("H3").Value = "L."
$cellI = $ws.Range("I3")
$cellI.NumberFormat = "@"
$cellI.Value = "2"
$cellI.NumberFormat = "General"
$cellI.Style = "Normal"

# Row 4 <- old row 5
$ws.Range("A4").Value = 87535080
$ws.Range("B4").Value = 96313
$ws.Range("E4").Value = 223609
$ws.Range("Q4").Value = 693933.2341120628
$ws.Range("R4").Value = 6685507.204237238
$ws.Range("F4").Value = "Vanlig skogsknipprot"
$ws.Range("G4").Value = "Epipactis helleborine subsp. helleborine"
$ws.Range("J4").Value = "stjälkar/strån/skott"
$ws.Range("AC4").Value = "Förekomst längs väg nr 1120 S om vägbanan i innerslänt."
$cellH = $ws.Range("H4")
$cellH.Value = "'"
$cellH.Style = "Normal"
$cellI = $ws.Range("I4")
$cellI.NumberFormat = "@"
$cellI.Value = "2"
$cellI.NumberFormat = "General"
$cellI.Style = "Normal"

# Row 5 <- old row 6
$ws.Range("A5").Value = 87534605
$ws.Range("B5").Value = 96356
$ws.Range("E5").Value = 219847
$ws.Range("Q5").Value = 693940.7630101059
$ws.Range("R5").Value = 6685506.11301556
$ws.Range("F5").Value = "Tvåblad"
$ws.Range("G5").Value = "Neottia ovata"
$ws.Range("J5").Value = "stjälkar/strån/skott"
$ws.Range("AC5").Value = "Förekomst längs väg nr 1120 S om vägbanan i innerslänt, dikesbotten och ytterslänt."
$ws.Range("H5").Value = "(L.) Buff. & Fingerh."
$cellI = $ws.Range("I5")
$cellI.NumberFormat = "@"
$cellI.Value = "25"
$cellI.NumberFormat = "General"
$cellI.Style = "Normal"

# Row 6 <- old row 7
$ws.Range("A6").Value = 87534816
$ws.Range("B6").Value = 96356
$ws.Range("E6").Value = 219847
$ws.Range("Q6").Value = 693928.1621041195
$ws.Range("R6").Value = 6685508.923384908
$ws.Range("F6").Value = "Tvåblad"
$ws.Range("G6").Value = "Neottia ovata"
$ws.Range("J6").Value = "stjälkar/strån/skott"
$ws.Range("AC6").Value = "Förekomst längs väg nr 1120 S om vägbanan i innerslänt."
$ws.Range("H6").Value = "(L.) Buff. & Fingerh."
$cellI = $ws.Range("I6")
$cellI.NumberFormat = "@"
$cellI.Value = "10"
$cellI.NumberFormat = "General"
$cellI.Style = "Normal"

# Row 7 <- old row 3
$ws.Range("A7").Value = 87534675
$ws.Range("B7").Value = 96356
$ws.Range("E7").Value = 219847
$ws.Range("Q7").Value = 693482.7895264921
$ws.Range("R7").Value = 6685726.906257272
$ws.Range("F7").Value = "Tvåblad"
$ws.Range("G7").Value = "Neottia ovata"
$ws.Range("J7").Value = "stjälkar/strån/skott"
$ws.Range("AC7").Value = "Förekomst längs väg nr 1120 S om vägbanan i innerslänt, dikesbotten och ytterslänt."
$ws.Range("H7").Value = "(L.) Buff. & Fingerh."
$cellI = $ws.Range("I7")
$cellI.NumberFormat = "@"
$cellI.Value = "16"
$cellI.NumberFormat = "General"
$cellI.Style = "Normal"

# AF4 gains an empty text cell (previously only on row 5)
$cellAF4 = $ws.Range("AF4")
$cellAF4.Value = "'"
$cellAF4.Style = "Normal"

# AF5 loses its (empty) cell entirely
$ws.Range("AF5").Value = ""

